$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 23
$ws1.Range("F4").Value = 87
$ws1.Range("F7").Value = 1706
$ws1.Range("F8").Value = 30
$ws1.Range("F11").Value = 1753
$ws1.Range("F12").Value = 129
$ws1.Range("F13").Value = 103
$ws1.Range("F14").Value = 418
$ws1.Range("F15").Value = 270
$ws1.Range("F16").Value = 198
$ws1.Range("F17").Value = 15
$ws1.Range("F18").Value = 31
$ws1.Range("F19").Value = 36
$ws1.Range("F21").Value = 522
$ws1.Range("F22").Value = 304
$ws1.Range("F23").Value = 165
$ws1.Range("F24").Value = 241
$ws1.Range("F25").Value = 256

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 23
$ws4.Range("F4").Value = 87
$ws4.Range("F7").Value = 1706
$ws4.Range("F9").Value = 30
$ws4.Range("F12").Value = 1753
$ws4.Range("F13").Value = 129
$ws4.Range("F14").Value = 103
$ws4.Range("F15").Value = 418
$ws4.Range("F16").Value = 270
$ws4.Range("F17").Value = 198
$ws4.Range("F18").Value = 15
$ws4.Range("F19").Value = 31
$ws4.Range("F20").Value = 36
$ws4.Range("F22").Value = 522
$ws4.Range("F23").Value = 304
$ws4.Range("F24").Value = 165
$ws4.Range("F25").Value = 241
$ws4.Range("F26").Value = 256
